$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap data for row pairs (2,3), (4,5), (11,12); column A (id) stays unchanged ---
# Row 2
$ws.Cells.Item(2, 2).Value = 6760228
$ws.Cells.Item(2, 3).Value = "Germany Landesliga"
$ws.Cells.Item(2, 4).Value = "Germany Landesliga"
$ws.Cells.Item(2, 5).Value = 45088.41666666666
$ws.Cells.Item(2, 6).Value = "SpVg Porz 1919"
$ws.Cells.Item(2, 7).Value = "Bonn Endenich 1908"
$ws.Cells.Item(2, 8).Value = 4
$ws.Cells.Item(2, 9).Value = 1
$ws.Cells.Item(2, 10).Value = "H"
$ws.Cells.Item(2, 11).Value = 2.25
$ws.Cells.Item(2, 12).Value = 3.5
$ws.Cells.Item(2, 13).Value = 2.625
$ws.Cells.Item(2, 14).Value = 1.75
$ws.Cells.Item(2, 15).Value = 3.6
$ws.Cells.Item(2, 16).Value = 3.8
$ws.Cells.Item(2, 17).Value = -0.5
$ws.Cells.Item(2, 18).Value = 1.8
$ws.Cells.Item(2, 19).Value = 2
$ws.Cells.Item(2, 20).Value = 2.75
$ws.Cells.Item(2, 21).Value = 1.8
$ws.Cells.Item(2, 22).Value = 2
$ws.Cells.Item(2, 23).Value = 0.75
$ws.Cells.Item(2, 24).Value = -1
$ws.Cells.Item(2, 25).Value = -1
$ws.Cells.Item(2, 26).Value = 0.8
$ws.Cells.Item(2, 27).Value = -1
$ws.Cells.Item(2, 28).Value = 0.8
$ws.Cells.Item(2, 29).Value = -1

# Row 3
$ws.Cells.Item(3, 2).Value = 6757276
$ws.Cells.Item(3, 3).Value = "Germany Landesliga"
$ws.Cells.Item(3, 4).Value = "Germany Landesliga"
$ws.Cells.Item(3, 5).Value = 45088.41666666666
$ws.Cells.Item(3, 6).Value = "BSC Rapid Chemnitz"
$ws.Cells.Item(3, 7).Value = "FV Dresden 06 Laubegast"
$ws.Cells.Item(3, 8).Value = 3
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = "H"
$ws.Cells.Item(3, 11).Value = 2.6
$ws.Cells.Item(3, 12).Value = 3.5
$ws.Cells.Item(3, 13).Value = 2.25
$ws.Cells.Item(3, 14).Value = 2.6
$ws.Cells.Item(3, 15).Value = 3.5
$ws.Cells.Item(3, 16).Value = 2.25
$ws.Cells.Item(3, 17).Value = 0.25
$ws.Cells.Item(3, 18).Value = 1.75
$ws.Cells.Item(3, 19).Value = 2.05
$ws.Cells.Item(3, 20).Value = 3.25
$ws.Cells.Item(3, 21).Value = 1.775
$ws.Cells.Item(3, 22).Value = 2.025
$ws.Cells.Item(3, 23).Value = 1.6
$ws.Cells.Item(3, 24).Value = -1
$ws.Cells.Item(3, 25).Value = -1
$ws.Cells.Item(3, 26).Value = 0.75
$ws.Cells.Item(3, 27).Value = -1
$ws.Cells.Item(3, 28).Value = -0.5
$ws.Cells.Item(3, 29).Value = 0.5125

# Row 4
$ws.Cells.Item(4, 2).Value = 6781315
$ws.Cells.Item(4, 3).Value = "Germany Landesliga"
$ws.Cells.Item(4, 4).Value = "Germany Landesliga"
$ws.Cells.Item(4, 5).Value = 45094.41666666666
$ws.Cells.Item(4, 6).Value = "SSV Markranstadt"
$ws.Cells.Item(4, 7).Value = "BSC Rapid Chemnitz"
$ws.Cells.Item(4, 8).Value = 2
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = "H"
$ws.Cells.Item(4, 11).Value = 1.25
$ws.Cells.Item(4, 12).Value = 4.75
$ws.Cells.Item(4, 13).Value = 10
$ws.Cells.Item(4, 14).Value = 1.222
$ws.Cells.Item(4, 15).Value = 5.25
$ws.Cells.Item(4, 16).Value = 8.5
$ws.Cells.Item(4, 17).Value = -2
$ws.Cells.Item(4, 18).Value = 1.925
$ws.Cells.Item(4, 19).Value = 1.875
$ws.Cells.Item(4, 20).Value = 3.5
$ws.Cells.Item(4, 21).Value = 1.775
$ws.Cells.Item(4, 22).Value = 1.925
$ws.Cells.Item(4, 23).Value = 0.222
$ws.Cells.Item(4, 24).Value = -1
$ws.Cells.Item(4, 25).Value = -1
$ws.Cells.Item(4, 26).Value = 0
$ws.Cells.Item(4, 27).Value = -0
$ws.Cells.Item(4, 28).Value = -1
$ws.Cells.Item(4, 29).Value = 0.925

# Row 5
$ws.Cells.Item(5, 2).Value = 6781316
$ws.Cells.Item(5, 3).Value = "Germany Landesliga"
$ws.Cells.Item(5, 4).Value = "Germany Landesliga"
$ws.Cells.Item(5, 5).Value = 45094.41666666666
$ws.Cells.Item(5, 6).Value = "SV Schott Jena"
$ws.Cells.Item(5, 7).Value = "SV 09 Arnstadt"
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 2
$ws.Cells.Item(5, 10).Value = "A"
$ws.Cells.Item(5, 11).Value = 5
$ws.Cells.Item(5, 12).Value = 4.5
$ws.Cells.Item(5, 13).Value = 1.45
$ws.Cells.Item(5, 14).Value = 6.5
$ws.Cells.Item(5, 15).Value = 4.333
$ws.Cells.Item(5, 16).Value = 1.363
$ws.Cells.Item(5, 17).Value = 1.5
$ws.Cells.Item(5, 18).Value = 1.825
$ws.Cells.Item(5, 19).Value = 1.975
$ws.Cells.Item(5, 20).Value = 3
$ws.Cells.Item(5, 21).Value = 1.825
$ws.Cells.Item(5, 22).Value = 1.975
$ws.Cells.Item(5, 23).Value = -1
$ws.Cells.Item(5, 24).Value = -1
$ws.Cells.Item(5, 25).Value = 0.363
$ws.Cells.Item(5, 26).Value = -1
$ws.Cells.Item(5, 27).Value = 0.9750000000000001
$ws.Cells.Item(5, 28).Value = -1
$ws.Cells.Item(5, 29).Value = 0.9750000000000001

# Row 11
$ws.Cells.Item(11, 2).Value = 7035046
$ws.Cells.Item(11, 3).Value = "Germany Landesliga"
$ws.Cells.Item(11, 4).Value = "Germany Landesliga"
$ws.Cells.Item(11, 5).Value = 45147.625
$ws.Cells.Item(11, 6).Value = "Cronenberger SC"
$ws.Cells.Item(11, 7).Value = "FC Viersen"
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 9).Value = 2
$ws.Cells.Item(11, 10).Value = "A"
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 3.6
$ws.Cells.Item(11, 13).Value = 3
$ws.Cells.Item(11, 14).Value = 2
$ws.Cells.Item(11, 15).Value = 3.6
$ws.Cells.Item(11, 16).Value = 3
$ws.Cells.Item(11, 17).Value = -0.25
$ws.Cells.Item(11, 18).Value = 1.8
$ws.Cells.Item(11, 19).Value = 2
$ws.Cells.Item(11, 20).Value = 2.75
$ws.Cells.Item(11, 21).Value = 1.8
$ws.Cells.Item(11, 22).Value = 2
$ws.Cells.Item(11, 23).Value = -1
$ws.Cells.Item(11, 24).Value = -1
$ws.Cells.Item(11, 25).Value = 2
$ws.Cells.Item(11, 26).Value = -1
$ws.Cells.Item(11, 27).Value = 1
$ws.Cells.Item(11, 28).Value = -1
$ws.Cells.Item(11, 29).Value = 1

# Row 12
$ws.Cells.Item(12, 2).Value = 7035048
$ws.Cells.Item(12, 3).Value = "Germany Landesliga"
$ws.Cells.Item(12, 4).Value = "Germany Landesliga"
$ws.Cells.Item(12, 5).Value = 45147.625
$ws.Cells.Item(12, 6).Value = "SG Unterrath"
$ws.Cells.Item(12, 7).Value = "TuRU Dsseldorf"
$ws.Cells.Item(12, 8).Value = 1
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = "H"
$ws.Cells.Item(12, 11).Value = 3.25
$ws.Cells.Item(12, 12).Value = 4
$ws.Cells.Item(12, 13).Value = 1.8
$ws.Cells.Item(12, 14).Value = 2.9
$ws.Cells.Item(12, 15).Value = 4
$ws.Cells.Item(12, 16).Value = 1.95
$ws.Cells.Item(12, 17).Value = 0.5
$ws.Cells.Item(12, 18).Value = 1.8
$ws.Cells.Item(12, 19).Value = 2
$ws.Cells.Item(12, 20).Value = 3
$ws.Cells.Item(12, 21).Value = 1.75
$ws.Cells.Item(12, 22).Value = 1.95
$ws.Cells.Item(12, 23).Value = 1.9
$ws.Cells.Item(12, 24).Value = -1
$ws.Cells.Item(12, 25).Value = -1
$ws.Cells.Item(12, 26).Value = 0.8
$ws.Cells.Item(12, 27).Value = -1
$ws.Cells.Item(12, 28).Value = -1
$ws.Cells.Item(12, 29).Value = 0.95

# --- New rows 80 and 81: first clone formatting from row 79, then set values ---
$ws.Range("A79:AC79").Copy()
$ws.Range("A80:AC80").PasteSpecial(-4122)
$ws.Range("A81:AC81").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 80
$ws.Cells.Item(80, 1).Value = 78
$ws.Cells.Item(80, 2).Value = 8059643
$ws.Cells.Item(80, 3).Value = "Germany Landesliga"
$ws.Cells.Item(80, 4).Value = "Germany Landesliga"
$ws.Cells.Item(80, 5).Value = 45389.4375
$ws.Cells.Item(80, 6).Value = "SC Victoria Mennrath"
$ws.Cells.Item(80, 7).Value = "SG Unterrath"
$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = "D"
$ws.Cells.Item(80, 11).Value = 2.2
$ws.Cells.Item(80, 12).Value = 3.5
$ws.Cells.Item(80, 13).Value = 2.7
$ws.Cells.Item(80, 14).Value = 2.2
$ws.Cells.Item(80, 15).Value = 3.6
$ws.Cells.Item(80, 16).Value = 2.625
$ws.Cells.Item(80, 17).Value = -0.25
$ws.Cells.Item(80, 18).Value = 2
$ws.Cells.Item(80, 19).Value = 1.8
$ws.Cells.Item(80, 20).Value = 3.25
$ws.Cells.Item(80, 21).Value = 1.8
$ws.Cells.Item(80, 22).Value = 2
$ws.Cells.Item(80, 23).Value = -1
$ws.Cells.Item(80, 24).Value = 2.6
$ws.Cells.Item(80, 25).Value = -1
$ws.Cells.Item(80, 26).Value = -0.5
$ws.Cells.Item(80, 27).Value = 0.4
$ws.Cells.Item(80, 28).Value = -1
$ws.Cells.Item(80, 29).Value = 1

# Row 81
$ws.Cells.Item(81, 1).Value = 79
$ws.Cells.Item(81, 2).Value = 8059644
$ws.Cells.Item(81, 3).Value = "Germany Landesliga"
$ws.Cells.Item(81, 4).Value = "Germany Landesliga"
$ws.Cells.Item(81, 5).Value = 45389.4375
$ws.Cells.Item(81, 6).Value = "ASV Suchteln"
$ws.Cells.Item(81, 7).Value = "FC Monheim"
$ws.Cells.Item(81, 8).Value = 2
$ws.Cells.Item(81, 9).Value = 4
$ws.Cells.Item(81, 10).Value = "A"
$ws.Cells.Item(81, 11).Value = 3.75
$ws.Cells.Item(81, 12).Value = 4
$ws.Cells.Item(81, 13).Value = 1.666
$ws.Cells.Item(81, 14).Value = 4.75
$ws.Cells.Item(81, 15).Value = 4.2
$ws.Cells.Item(81, 16).Value = 1.5
$ws.Cells.Item(81, 17).Value = 1
$ws.Cells.Item(81, 18).Value = 1.975
$ws.Cells.Item(81, 19).Value = 1.825
$ws.Cells.Item(81, 20).Value = 3.5
$ws.Cells.Item(81, 21).Value = 1.9
$ws.Cells.Item(81, 22).Value = 1.9
$ws.Cells.Item(81, 23).Value = -1
$ws.Cells.Item(81, 24).Value = -1
$ws.Cells.Item(81, 25).Value = 0.5
$ws.Cells.Item(81, 26).Value = -1
$ws.Cells.Item(81, 27).Value = 0.825
$ws.Cells.Item(81, 28).Value = 0.8999999999999999
$ws.Cells.Item(81, 29).Value = -1

